# Fix #354, #360, #372
# - plotConfiguration: add "quantiles" and "foldDistance" parameter columns,
#   and set a "foldDistance" value of "2, 3" for the observedVsSimulated plot (P2).
# - plotGrids: add a "tagPrefix" parameter column.

$wb = $excel.ActiveWorkbook

# --- plotConfiguration sheet -------------------------------------------------
$wsPlotConfig = $wb.Worksheets.Item("plotConfiguration")

$wsPlotConfig.Range("K1").Value = "quantiles"
$wsPlotConfig.Range("L1").Value = "foldDistance"

# Row 3 corresponds to plotID P2 / observedVsSimulated
$wsPlotConfig.Range("L3").Value = "2, 3"

# --- plotGrids sheet ---------------------------------------------------------
$wsPlotGrids = $wb.Worksheets.Item("plotGrids")

$wsPlotGrids.Range("D1").Value = "tagPrefix"

# --- exportConfiguration sheet ----------------------------------------------
# Update its remembered selection before leaving it.
$wsExportConfig = $wb.Worksheets.Item("exportConfiguration")
$wsExportConfig.Activate()
$wsExportConfig.Range("B3").Select()

# Update plotConfiguration's remembered selection to the newly-filled cell.
$wsPlotConfig.Activate()
$wsPlotConfig.Range("L3").Select()

# Make plotGrids the active sheet/selection, matching the saved view state.
$wsPlotGrids.Activate()
$wsPlotGrids.Range("D3").Select()
